$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prepay Loan" (index 2): update payoff value, move selection to B5
# ---------------------------------------------------------------------
$wsPrepay = $wb.Worksheets.Item(2)
$wsPrepay.Activate() | Out-Null
$wsPrepay.Range("B5").Value = 10015.34
$wsPrepay.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Summary" (index 3): refresh overdue figures, move selection
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item(3)
$wsSummary.Activate() | Out-Null

$wsSummary.Range("B2").Value = 9981.3700000000008
$wsSummary.Range("B2").NumberFormat = "#,##0.00"

$wsSummary.Range("E2").Value = 18.63
$wsSummary.Range("F2").Value = 18.63

# A3/B3 used to hold formulas; replace with plain overdue numbers and
# restore the plain (non-formula) number format used elsewhere in the row
$wsSummary.Range("C2").Copy() | Out-Null
$wsSummary.Range("A3:B3").PasteSpecial(-4122) | Out-Null
$wsSummary.Range("A3").Value = 34.19
$wsSummary.Range("B3").Value = 33.97

$wsSummary.Range("E3").Value = 0.22
$wsSummary.Range("F3").Value = 0.09

$wsSummary.Range("A7:XFD14").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Repayment schedule" (index 4): drop the now-unused trailing
# blank rows 5-14 and move the selection
# ---------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item(4)
$wsRepay.Activate() | Out-Null
$wsRepay.Range("A5:R14").EntireRow.Delete() | Out-Null
$wsRepay.Range("A5:XFD16").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Sheet1" (index 6) becomes the active tab
# ---------------------------------------------------------------------
$wsSheet1 = $wb.Worksheets.Item(6)
$wsSheet1.Activate() | Out-Null
